# Apply "Baseline Wandering Cleared" score updates to the scores table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - 1dAVb
$ws.Range("B4").Value = 1
$ws.Range("F4").Value = 0.964
$ws.Range("J4").Value = 1
$ws.Range("N4").Value = 0.982

# Row 5 - RBBB
$ws.Range("B5").Value = 0.872
$ws.Range("J5").Value = 0.994
$ws.Range("N5").Value = 0.9320000000000001

# Row 6 - LBBB
$ws.Range("F6").Value = 1
$ws.Range("N6").Value = 1

# Row 7 - SB
$ws.Range("B7").Value = 0.889
$ws.Range("F7").Value = 1
$ws.Range("J7").Value = 0.998
$ws.Range("N7").Value = 0.9409999999999999

# Row 8 - AF
$ws.Range("F8").Value = 0.385
$ws.Range("N8").Value = 0.556

# Row 9 - ST
$ws.Range("B9").Value = 1
$ws.Range("F9").Value = 0.865
$ws.Range("J9").Value = 1
$ws.Range("N9").Value = 0.928
